# New crime data collected - weekly CompStat refresh.
# Updates the "Volume/Number" + reporting-week header text and refreshes the
# crime-stat grid (rows 16-27) with the next week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Plain numeric write: leaves the cell's existing style/number-format alone.
function SetNumber {
    param($addr, $value)
    $ws.Range($addr).Value = $value
}

# Numeric write into a cell that currently holds a text placeholder
# (shared string "0" / "***.*", style 14). Setting NumberFormat picks the
# matching pre-existing numeric style (16 for integers, 15 for percentages)
# instead of minting a new one.
function SetNumberWithFormat {
    param($addr, $value, $numFmt)
    $cell = $ws.Range($addr)
    $cell.Value = $value
    $cell.NumberFormat = $numFmt
}

# Text write into a cell that currently holds a number. Forcing Text format
# first stops the engine from re-coercing a numeric-looking string ("0")
# back into a number; the subsequent format-only paste from a known-good
# style-14 cell restores the normal (General) right-aligned look so we don't
# leave the cell on a bespoke "@" format.
function SetText {
    param($addr, $value)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $styleSource = $ws.Range("C14")
    $styleSource.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# Header text: Volume/Number + reporting week
# ---------------------------------------------------------------------------

# "Volume 30   Number  19" -> "Volume 30   Number  20"
$ws.Range("A8").Characters(21, 2).Text = "20"

# "Report Covering the Week  5/8/2023  Through  5/14/2023"
#   -> "Report Covering the Week  5/15/2023  Through  5/21/2023"
# Replace the later (Through) date first so the earlier date's character
# offsets stay valid.
$ws.Range("C9").Characters(46, 9).Text = "5/21/2023"
$ws.Range("C9").Characters(27, 8).Text = "5/15/2023"

# ---------------------------------------------------------------------------
# Crime grid refresh (rows 16-27)
# ---------------------------------------------------------------------------

# Row 16 - Robbery
SetNumberWithFormat "D16" 2 "#,##0"
SetNumberWithFormat "E16" -50 '#,##0.0;"-"#,##0.0'
SetNumber "F16" 2
SetNumber "G16" 3
SetNumber "H16" -33.333333333333
SetNumber "I16" 14
SetNumber "J16" 15
SetNumber "K16" -6.666666666666
SetNumber "L16" 75
SetNumber "M16" -41.666666666666
SetNumber "N16" -86.274509803921

# Row 17 - Fel. Assault
SetNumber "C17" 2
SetNumber "D17" 3
SetNumber "E17" -33.333333333333
SetNumber "F17" 10
SetNumber "G17" 10
SetNumber "I17" 33
SetNumber "J17" 36
SetNumber "K17" -8.333333333333
SetNumber "L17" -2.941176470588
SetNumber "M17" 6.451612903225
SetNumber "N17" -57.142857142857

# Row 18 - Burglary
SetText "D18" "0"
SetText "E18" "***.*"
SetNumber "F18" 3
SetNumber "G18" 3
SetNumber "H18" 0
SetNumber "I18" 13
SetNumber "K18" -38.095238095238
SetNumber "L18" -13.333333333333
SetNumber "M18" -31.578947368421
SetNumber "N18" -93.467336683417

# Row 19 - Gr. Larceny
SetNumberWithFormat "D19" 5 "#,##0"
SetNumberWithFormat "E19" -60 '#,##0.0;"-"#,##0.0'
SetNumber "F19" 7
SetNumber "G19" 8
SetNumber "H19" -12.5
SetNumber "I19" 62
SetNumber "J19" 52
SetNumber "K19" 19.230769230769
SetNumber "L19" 55
SetNumber "M19" 77.142857142857
SetNumber "N19" 14.814814814814

# Row 20 - G.L.A.
SetText "D20" "0"
SetText "E20" "***.*"
SetNumber "F20" 3
SetNumber "G20" 2
SetNumber "H20" 50
SetNumber "L20" 37.5
SetNumber "N20" -92.028985507246

# Row 21 - TOTAL (bold styles 18/19; values only, styles unaffected)
SetNumber "C21" 6
SetNumber "D21" 10
SetNumber "E21" -40
SetNumber "F21" 25
SetNumber "G21" 26
SetNumber "H21" -3.846153846153
SetNumber "I21" 134
SetNumber "J21" 142
SetNumber "K21" -5.633802816901
SetNumber "L21" 21.818181818181
SetNumber "M21" 12.605042016806
SetNumber "N21" -76.856649395509

# Row 23 - Housing
SetText "D23" "0"
SetText "E23" "***.*"
SetText "F23" "0"
SetNumber "G23" 3
SetNumber "H23" -100
SetNumber "L23" -22.222222222222

# Row 24 - Petit Larceny
SetNumber "C24" 6
SetNumber "D24" 6
SetNumber "E24" 0
SetNumber "F24" 31
SetNumber "G24" 34
SetNumber "H24" -8.823529411764
SetNumber "I24" 177
SetNumber "J24" 186
SetNumber "K24" -4.838709677419
SetNumber "L24" 33.082706766917
SetNumber "M24" 90.322580645161

# Row 25 - Misd. Assault
SetNumber "C25" 5
SetNumber "D25" 7
SetNumber "E25" -28.571428571428
SetNumber "G25" 19
SetNumber "H25" 10.526315789473
SetNumber "I25" 77
SetNumber "J25" 68
SetNumber "K25" 13.235294117647
SetNumber "L25" 28.333333333333
SetNumber "M25" -42.105263157894

# Row 27 - Other Sex Crimes
SetText "D27" "0"
SetText "E27" "***.*"
SetNumber "F27" 1
SetNumber "H27" -50
